$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,3).Value = 0.2804546356201172
$ws.Cells.Item(2,4).Value = 0.4303635954856872
$ws.Cells.Item(2,5).Value = -0.691750168800354
$ws.Cells.Item(2,6).Value = 0.1050096067542932
$ws.Cells.Item(2,7).Value = -1.756468223065746
$ws.Cells.Item(2,8).Value = 0.4945203567645989
$ws.Cells.Item(3,3).Value = 0.1987819671630859
$ws.Cells.Item(3,4).Value = 0.2879692316055298
$ws.Cells.Item(3,5).Value = -0.9282988905906676
$ws.Cells.Item(3,6).Value = -0.0286234012063665
$ws.Cells.Item(3,7).Value = -0.7998002785809195
$ws.Cells.Item(3,8).Value = 0.0811297598541999
$ws.Cells.Item(4,3).Value = 0.2130537033081054
$ws.Cells.Item(4,4).Value = 0.492556095123291
$ws.Cells.Item(4,5).Value = -0.9500083923339844
$ws.Cells.Item(4,6).Value = 0.03028146071093407
$ws.Cells.Item(4,7).Value = 0.6503314929349078
$ws.Cells.Item(4,8).Value = -0.2584614434412545
$ws.Cells.Item(5,3).Value = 5.370540142059326
$ws.Cells.Item(5,4).Value = -2.71519660949707
$ws.Cells.Item(5,5).Value = -10.77233219146728
$ws.Cells.Item(5,6).Value = 0.6556547600395809
$ws.Cells.Item(5,7).Value = 3.617322595752007
$ws.Cells.Item(5,8).Value = -1.325163067603595
$ws.Cells.Item(6,3).Value = 1.807073593139648
$ws.Cells.Item(6,4).Value = 1.692888975143433
$ws.Cells.Item(6,5).Value = 0.979395866394043
$ws.Cells.Item(6,6).Value = 0.9173419645854408
$ws.Cells.Item(6,7).Value = 4.483878631981052
$ws.Cells.Item(6,8).Value = -1.026156997194096
$ws.Cells.Item(7,3).Value = -1.115023136138916
$ws.Cells.Item(7,4).Value = -0.1916569471359253
$ws.Cells.Item(7,5).Value = -1.226519584655761
$ws.Cells.Item(7,6).Value = 0.6556111379545562
$ws.Cells.Item(7,7).Value = 2.95567157317181
$ws.Cells.Item(7,8).Value = 0.3850134197546516
$ws.Cells.Item(8,3).Value = -0.6978027820587158
$ws.Cells.Item(8,4).Value = -0.2956833839416504
$ws.Cells.Item(8,5).Value = 0.6510324478149414
$ws.Cells.Item(8,6).Value = 0.01345462458474067
$ws.Cells.Item(8,7).Value = 2.573593986277679
$ws.Cells.Item(8,8).Value = -0.1254206044333366
$ws.Cells.Item(9,3).Value = -3.345348596572876
$ws.Cells.Item(9,4).Value = 1.023304224014282
$ws.Cells.Item(9,5).Value = 2.694005012512207
$ws.Cells.Item(9,6).Value = -0.0735905529284972
$ws.Cells.Item(9,7).Value = 1.541839412736652
$ws.Cells.Item(9,8).Value = -0.5845232478209913
$ws.Cells.Item(10,3).Value = -9.996002197265623
$ws.Cells.Item(10,4).Value = -0.0120857954025268
$ws.Cells.Item(10,5).Value = 1.461801767349243
$ws.Cells.Item(10,6).Value = 0.5347626531610679
$ws.Cells.Item(10,7).Value = -1.351723280640276
$ws.Cells.Item(10,8).Value = 0.006906517306153076
$ws.Cells.Item(11,3).Value = 4.229874610900879
$ws.Cells.Item(11,4).Value = -2.616225481033325
$ws.Cells.Item(11,5).Value = -11.63714218139648
$ws.Cells.Item(11,6).Value = 0.2076723618166787
$ws.Cells.Item(11,7).Value = -3.982405117579868
$ws.Cells.Item(11,8).Value = -0.1024508295314652
$ws.Cells.Item(12,3).Value = -5.585046291351318
$ws.Cells.Item(12,4).Value = -1.405134916305542
$ws.Cells.Item(12,5).Value = 1.352597236633301
$ws.Cells.Item(12,6).Value = -0.6092913156869469
$ws.Cells.Item(12,7).Value = -4.47834036301593
$ws.Cells.Item(12,8).Value = 0.3969252789203029
$ws.Cells.Item(13,3).Value = -3.109971046447754
$ws.Cells.Item(13,4).Value = -0.2804408073425293
$ws.Cells.Item(13,5).Value = -1.378412246704102
$ws.Cells.Item(13,6).Value = -0.4407143246154432
$ws.Cells.Item(13,7).Value = -3.636212723595754
$ws.Cells.Item(13,8).Value = 0.6680403655889118
$ws.Cells.Item(14,3).Value = -6.216216564178467
$ws.Cells.Item(14,4).Value = 2.080293655395508
$ws.Cells.Item(14,5).Value = -6.55400562286377
$ws.Cells.Item(14,6).Value = 0.2482107579708095
$ws.Cells.Item(14,7).Value = -2.94579162889598
$ws.Cells.Item(14,8).Value = 0.002957710806204883
$ws.Cells.Item(15,3).Value = -2.478143692016602
$ws.Cells.Item(15,4).Value = -0.5579037666320801
$ws.Cells.Item(15,5).Value = -6.963788509368896
$ws.Cells.Item(15,6).Value = 0.1622283559064473
$ws.Cells.Item(15,7).Value = 0.07419828006199625
$ws.Cells.Item(15,8).Value = -0.1648120685499541
$ws.Cells.Item(16,3).Value = 1.977913856506348
$ws.Cells.Item(16,4).Value = 1.696393609046936
$ws.Cells.Item(16,5).Value = 4.733741760253906
$ws.Cells.Item(16,6).Value = 0.2419649730531533
$ws.Cells.Item(16,7).Value = 4.375038618944129
$ws.Cells.Item(16,8).Value = -0.3818312956362353
$ws.Cells.Item(17,3).Value = -1.786364078521728
$ws.Cells.Item(17,4).Value = -0.5252545475959778
$ws.Cells.Item(17,5).Value = -2.082462787628174
$ws.Cells.Item(17,6).Value = 0.5893073264433398
$ws.Cells.Item(17,7).Value = 6.328164937544838
$ws.Cells.Item(17,8).Value = 0.07646724885823974
$ws.Cells.Item(18,3).Value = -3.566806316375732
$ws.Cells.Item(18,4).Value = 0.4948284029960632
$ws.Cells.Item(18,5).Value = 1.737120628356934
$ws.Cells.Item(18,6).Value = -0.175296502453943
$ws.Cells.Item(18,7).Value = 2.29715876919883
$ws.Cells.Item(18,8).Value = 1.448470601013731
$ws.Cells.Item(19,3).Value = -6.714832305908203
$ws.Cells.Item(19,4).Value = 2.814615249633789
$ws.Cells.Item(19,5).Value = 5.518161296844482
$ws.Cells.Item(19,6).Value = 0.8341801327710269
$ws.Cells.Item(19,7).Value = 0.753240480714915
$ws.Cells.Item(19,8).Value = -0.1038439748238527
$ws.Cells.Item(20,3).Value = -4.827847480773926
$ws.Cells.Item(20,4).Value = -7.360194683074951
$ws.Cells.Item(20,5).Value = 7.510563373565674
$ws.Cells.Item(20,6).Value = -0.1397884144466732
$ws.Cells.Item(20,7).Value = 0.2036113617371567
$ws.Cells.Item(20,8).Value = -0.8063265510967792
$ws.Cells.Item(21,3).Value = -1.930130958557129
$ws.Cells.Item(21,4).Value = 4.015055179595947
$ws.Cells.Item(21,5).Value = -5.047464370727539
$ws.Cells.Item(21,6).Value = -0.7627712436476506
$ws.Cells.Item(21,7).Value = -1.678813515877241
$ws.Cells.Item(21,8).Value = -1.093177660387389
$ws.Cells.Item(22,1).Value = 2000
$ws.Cells.Item(22,2).Value = "struggle"
$ws.Cells.Item(22,3).Value = -4.561958312988281
$ws.Cells.Item(22,4).Value = 0.1262733936309814
$ws.Cells.Item(22,5).Value = -1.993290901184082
$ws.Cells.Item(22,6).Value = -0.0004051567948546442
$ws.Cells.Item(22,7).Value = -2.697911312993704
$ws.Cells.Item(22,8).Value = -0.4782794093599184
$ws.Cells.Item(23,1).Value = 2100
$ws.Cells.Item(23,2).Value = "struggle"
$ws.Cells.Item(23,3).Value = -2.729169845581055
$ws.Cells.Item(23,4).Value = 3.059413433074951
$ws.Cells.Item(23,5).Value = -4.533473014831543
$ws.Cells.Item(23,6).Value = -0.8633800915309378
$ws.Cells.Item(23,7).Value = -0.3154059344408438
$ws.Cells.Item(23,8).Value = 0.4824308418497782
$ws.Cells.Item(24,1).Value = 2200
$ws.Cells.Item(24,2).Value = "struggle"
$ws.Cells.Item(24,3).Value = -0.4514303207397461
$ws.Cells.Item(24,4).Value = -0.07753515243530271
$ws.Cells.Item(24,5).Value = -1.056098580360413
$ws.Cells.Item(24,6).Value = -0.4081483519807154
$ws.Cells.Item(24,7).Value = -0.6726997543354425
$ws.Cells.Item(24,8).Value = -0.2190668820118418
$ws.Cells.Item(25,1).Value = 2300
$ws.Cells.Item(25,2).Value = "struggle"
$ws.Cells.Item(25,3).Value = 1.037992477416992
$ws.Cells.Item(25,4).Value = -1.273390769958496
$ws.Cells.Item(25,5).Value = 0.4362349510192871
$ws.Cells.Item(25,6).Value = 0.2211332225373814
$ws.Cells.Item(25,7).Value = 0.241335413285664
$ws.Cells.Item(25,8).Value = 0.08368853798934378
$ws.Cells.Item(26,1).Value = 2400
$ws.Cells.Item(26,2).Value = "struggle"
$ws.Cells.Item(26,3).Value = 0.0754270553588867
$ws.Cells.Item(26,4).Value = 1.646718859672546
$ws.Cells.Item(26,5).Value = 1.695090532302856
$ws.Cells.Item(26,6).Value = 0.06768137718341787
$ws.Cells.Item(26,7).Value = 0.3379019900244107
$ws.Cells.Item(26,8).Value = 0.1505034766635118
$ws.Cells.Item(27,1).Value = 2500
$ws.Cells.Item(27,2).Value = "struggle"
$ws.Cells.Item(27,3).Value = -0.2560558319091797
$ws.Cells.Item(27,4).Value = 0.3026316165924072
$ws.Cells.Item(27,5).Value = -0.4233262538909912
$ws.Cells.Item(27,6).Value = 0.07254024853511698
$ws.Cells.Item(27,7).Value = 0.5556785336562575
$ws.Cells.Item(27,8).Value = -0.05807583201296457
$ws.Cells.Item(28,1).Value = 2600
$ws.Cells.Item(28,2).Value = "struggle"
$ws.Cells.Item(28,3).Value = 0.6335611343383789
$ws.Cells.Item(28,4).Value = 0.8106564879417419
$ws.Cells.Item(28,5).Value = -1.443797469139099
$ws.Cells.Item(28,6).Value = 0.1816357883567719
$ws.Cells.Item(28,7).Value = 0.1322741392923868
$ws.Cells.Item(28,8).Value = -0.08515337003128903
$ws.Cells.Item(29,1).Value = 2700
$ws.Cells.Item(29,2).Value = "struggle"
$ws.Cells.Item(29,3).Value = 0.09285736083984369
$ws.Cells.Item(29,4).Value = 0.7357764840126038
$ws.Cells.Item(29,5).Value = -1.646607518196106
$ws.Cells.Item(29,6).Value = -0.02734556931013958
$ws.Cells.Item(29,7).Value = -0.1169588795425942
$ws.Cells.Item(29,8).Value = 0.04497027853313797
$ws.Cells.Item(30,1).Value = 2800
$ws.Cells.Item(30,2).Value = "struggle"
$ws.Cells.Item(30,3).Value = 0.0882749557495117
$ws.Cells.Item(30,4).Value = 0.1726978719234466
$ws.Cells.Item(30,5).Value = -0.9354652166366576
$ws.Cells.Item(30,6).Value = -0.02540700723017953
$ws.Cells.Item(30,7).Value = -0.06986615411481072
$ws.Cells.Item(30,8).Value = -0.074921377335808
$ws.Cells.Item(31,1).Value = 2900
$ws.Cells.Item(31,2).Value = "struggle"
$ws.Cells.Item(31,3).Value = 0.2656211853027344
$ws.Cells.Item(31,4).Value = 0.4902379512786865
$ws.Cells.Item(31,5).Value = -0.8409426212310791
$ws.Cells.Item(31,6).Value = 0.02237761537639455
$ws.Cells.Item(31,7).Value = -0.07008743807863513
$ws.Cells.Item(31,8).Value = -0.003453258577050004
